$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = -8
$ws.Range("F7").Value = 0
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = 3
$ws.Range("F13").Value = -4
$ws.Range("F15").Value = -4
$ws.Range("F17").Value = 2
